# Update cryptos worksheet with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '37.071.20'
    'E2' = '  -0.37%  '
    'D3' = '2.084.14'
    'E3' = '  +8.39%  '
    'E4' = '  -0.03%  '
    'D5' = '250.77'
    'E5' = '  +0.45%  '
    'D6' = '0.654'
    'E6' = '  -5.03%  '
    'E7' = '  +0.17%  '
    'D8' = '49.57'
    'E8' = '  +4.12%  '
    'D9' = '60.42'
    'E9' = '  +3.36%  '
    'D10' = '0.370'
    'E10' = '  -1.59%  '
    'D11' = '0.0740'
    'E11' = '  -3.16%  '
    'E12' = '  +5.06%  '
    'D13' = '15.04'
    'E13' = '  -3.86%  '
    'D14' = '2.386.37'
    'E14' = '  +8.30%  '
    'D15' = '0.827'
    'E15' = '  +0.21%  '
    'D16' = '2.091.69'
    'E16' = '  +8.77%  '
    'D17' = '5.08'
    'E17' = '  -1.02%  '
    'D18' = '36.930.61'
    'E18' = '  -0.82%  '
    'D19' = '71.92'
    'E19' = '  -4.02%  '
    'D20' = '0.0₃0820'
    'E20' = '  -4.45%  '
    'D21' = '13.20'
    'E21' = '  -3.34%  '
    'D22' = '238.96'
    'E22' = '  -4.79%  '
    'D23' = '5.17'
    'E23' = '  -0.16%  '
    'D24' = '0.999'
    'E24' = '  -0.19%  '
    'D25' = '2.45'
    'E25' = '  -2.73%  '
    'D26' = '168.29'
    'E26' = '  -0.01%  '
    'D27' = '9.26'
    'E27' = '  +5.14%  '
    'D28' = '20.72'
    'E28' = '  +10.57%  '
    'D29' = '2.01'
    'E29' = '  -4.84%  '
    'E30' = '  -5.07%  '
    'D31' = '23.57'
    'E31' = '  +23.01%  '
    'E32' = '  +25.92%  '
    'D33' = '4.45'
    'E33' = '  -2.28%  '
    'D34' = '0.0603'
    'E34' = '  -1.27%  '
    'D35' = '0.0915'
    'E35' = '  -0.26%  '
    'E36' = '  -0.04%  '
    'D37' = '2.27'
    'E37' = '  +16.72%  '
    'D38' = '1.83'
    'E38' = '  -3.12%  '
    'E39' = '  -5.27%  '
    'D40' = '1.31'
    'E40' = '  -9.88%  '
    'D41' = '17.58'
    'E41' = '  -0.32%  '
    'D42' = '0.0223'
    'E42' = '  -1.77%  '
    'D43' = '1.14'
    'E43' = '  +4.46%  '
    'D44' = '97.05'
    'E44' = '  -7.87%  '
    'D45' = '2.80'
    'E45' = '  -4.37%  '
    'D46' = '0.0867'
    'E46' = '  +3.71%  '
    'D47' = '2.96'
    'E47' = '  +5.26%  '
    'D48' = '1.303.90'
    'E48' = '  -3.10%  '
    'D49' = '6.85'
    'E49' = '  +6.79%  '
    'D50' = '2.259.16'
    'E50' = '  +6.98%  '
    'B51' = 'RenderToken'
    'C51' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D51' = '2.25'
    'E51' = '  -6.50%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
